# Add stage3 evidence / B1-B2
#
# The "Info" sheet is currently the active/selected sheet (cell D23
# selected). We move the selection on Info to B2 (no longer the active
# sheet), fill in the two new TxHash rows on sheets "B1" and "B2" (which
# appends four new shared strings, in this exact order so the shared-string
# table indices line up with the target), leave the selection on each of
# those sheets parked on A4 (the row below the newly written data), and
# finish with "B2" as the active sheet/tab.

$wb = $excel.ActiveWorkbook

# Info sheet: keep it unselected going forward, but its remembered
# selection moves from D23 to B2.
$wsInfo = $wb.Worksheets.Item("Info")
$wsInfo.Range("B2").Select()

# B1: write the stage-3 evidence TxHash values into A2/A3.
$wsB1 = $wb.Worksheets.Item("B1")
$wsB1.Range("A2").Value = "A8714DADFE7C5F64101B198DBFE466FCDDD277C15CBDAE7A36548C325EAD6CD9"
$wsB1.Range("A3").Value = "B0CA2688210CB1120BBFB2A020C3FCA4FC715E8707C1457C8349008441919B6F"
$wsB1.Range("A4").Select()

# B2: write the stage-3 evidence TxHash values into A2/A3, then make this
# the active sheet (tabSelected) with A4 as the resting selection.
$wsB2 = $wb.Worksheets.Item("B2")
$wsB2.Range("A2").Value = "7B762579D847876E32D782EB4456C813AB8EFE1F224B34FDB40F8A4F00FA56EC"
$wsB2.Range("A3").Value = "1BC72E7F6FAAA96374C241AD88A674C26FEE6ECDE15A454319AB75B764F783FA"
$wsB2.Activate()
$wsB2.Range("A4").Select()
